$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dSF column (F) values as per repulled data / mean calculation
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F14").Value = 0
